$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 2641.48
$ws.Range("I62").Value = 2329.5833
$ws.Range("J62").Value = 2929.3845
$ws.Range("K62").Value = 2329.5833
$ws.Range("L62").Value = 2929.3845
$ws.Range("M62").Value = -1705.5833
$ws.Range("N62").Value = -4177.3845
$ws.Range("H65").Value = 2641.48
$ws.Range("I65").Value = 2329.5833
$ws.Range("J65").Value = 2929.3845
$ws.Range("K65").Value = 11647.9165
$ws.Range("L65").Value = 14646.9225
$ws.Range("M65").Value = -8527.916499999999
$ws.Range("N65").Value = -20886.9225
$ws.Range("H74").Value = 6254310
$ws.Range("I74").Value = 3912.4443
$ws.Range("J74").Value = 11368271
$ws.Range("K74").Value = 3912.4443
$ws.Range("L74").Value = 11368271
$ws.Range("M74").Value = -2976.4443
$ws.Range("N74").Value = -11370143
$ws.Range("H77").Value = 6254310
$ws.Range("I77").Value = 3912.4443
$ws.Range("J77").Value = 11368271
$ws.Range("K77").Value = 19562.2215
$ws.Range("L77").Value = 56841355
$ws.Range("M77").Value = -14882.2215
$ws.Range("N77").Value = -56850715
$ws.Range("H86").Value = 5859.273
$ws.Range("I86").Value = 1265.6923
$ws.Range("K86").Value = 1265.6923
$ws.Range("M86").Value = -142.6922999999999
$ws.Range("H89").Value = 5859.273
$ws.Range("I89").Value = 1265.6923
$ws.Range("K89").Value = 6328.461499999999
$ws.Range("M89").Value = -712.4614999999994
$ws.Range("H96").Value = 1742
$ws.Range("I96").Value = 1525.5
$ws.Range("J96").Value = 1796.125
$ws.Range("K96").Value = 4576.5
$ws.Range("L96").Value = 5388.375
$ws.Range("M96").Value = -3203.5
$ws.Range("N96").Value = -8134.375
$ws.Range("H132").Value = 2874.6316
$ws.Range("I132").Value = 2996.258
$ws.Range("J132").Value = 2336
$ws.Range("K132").Value = 8988.773999999999
$ws.Range("L132").Value = 7008
$ws.Range("M132").Value = -6458.773999999999
$ws.Range("N132").Value = -12068

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6612.769
$ws.Range("I32").Value = 5227.3467
$ws.Range("J32").Value = 13106.9375
$ws.Range("K32").Value = 5227.3467
$ws.Range("L32").Value = 13106.9375
$ws.Range("M32").Value = -4940.3467
$ws.Range("N32").Value = -13680.9375
$ws.Range("H61").Value = 2598.3674
$ws.Range("I61").Value = 2746.2354
$ws.Range("J61").Value = 2263.2
$ws.Range("K61").Value = 2746.2354
$ws.Range("L61").Value = 2263.2
$ws.Range("M61").Value = -2534.2354
$ws.Range("N61").Value = -2687.2
$ws.Range("H74").Value = 26317168
$ws.Range("I74").Value = 32258678
$ws.Range("K74").Value = 32258678
$ws.Range("M74").Value = -32257804
$ws.Range("H77").Value = 26317168
$ws.Range("I77").Value = 32258678
$ws.Range("K77").Value = 161293390
$ws.Range("M77").Value = -161289022
$ws.Range("H102").Value = 1483.8
$ws.Range("I102").Value = 1404.3846
$ws.Range("K102").Value = 1404.3846
$ws.Range("M102").Value = 217.6153999999999
$ws.Range("H136").Value = 2598.3674
$ws.Range("I136").Value = 2746.2354
$ws.Range("J136").Value = 2263.2
$ws.Range("K136").Value = 8238.706200000001
$ws.Range("L136").Value = 6789.599999999999
$ws.Range("M136").Value = -5688.706200000001
$ws.Range("N136").Value = -11889.6

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1800.1
$ws.Range("I99").Value = 2166.6667
$ws.Range("J99").Value = 1643
$ws.Range("K99").Value = 2166.6667
$ws.Range("L99").Value = 1643
$ws.Range("M99").Value = -668.6667000000002
$ws.Range("N99").Value = -4639
$ws.Range("H134").Value = 3168.6904
$ws.Range("I134").Value = 3380.8157
$ws.Range("J134").Value = 1153.5
$ws.Range("K134").Value = 10142.4471
$ws.Range("L134").Value = 3460.5
$ws.Range("M134").Value = -7607.447100000001
$ws.Range("N134").Value = -8530.5

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4530.4316
$ws.Range("I31").Value = 2345.9412
$ws.Range("K31").Value = 2345.9412
$ws.Range("M31").Value = -2050.9412
$ws.Range("H34").Value = 4530.4316
$ws.Range("I34").Value = 2345.9412
$ws.Range("K34").Value = 2345.9412
$ws.Range("M34").Value = -2143.9412
$ws.Range("H47").Value = 15000
$ws.Range("J47").Value = 0
$ws.Range("L47").Value = 0
$ws.Range("N47").ClearContents()
$ws.Range("H99").Value = 3793.5
$ws.Range("I99").Value = 2853.625
$ws.Range("J99").Value = 6299.8335
$ws.Range("K99").Value = 2853.625
$ws.Range("L99").Value = 6299.8335
$ws.Range("M99").Value = -1355.625
$ws.Range("N99").Value = -9295.833500000001
$ws.Range("H126").Value = 3793.5
$ws.Range("I126").Value = 2853.625
$ws.Range("J126").Value = 6299.8335
$ws.Range("K126").Value = 8560.875
$ws.Range("L126").Value = 18899.5005
$ws.Range("M126").Value = -6090.875
$ws.Range("N126").Value = -23839.5005
$ws.Range("H132").Value = 3102.48
$ws.Range("I132").Value = 2216.25
$ws.Range("J132").Value = 4678
$ws.Range("K132").Value = 6648.75
$ws.Range("L132").Value = 14034
$ws.Range("M132").Value = -4118.75
$ws.Range("N132").Value = -19094

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 722.76
$ws.Range("I131").Value = 329.25
$ws.Range("J131").Value = 739.15625
$ws.Range("K131").Value = 987.75
$ws.Range("L131").Value = 2217.46875
$ws.Range("M131").Value = 4052.25
$ws.Range("N131").Value = -12297.46875
$ws.Range("H140").Value = 4497.875
$ws.Range("I140").Value = 3161.6667
$ws.Range("J140").Value = 5299.6
$ws.Range("K140").Value = 9485.000100000001
$ws.Range("L140").Value = 15898.8
$ws.Range("M140").Value = -4305.000100000001
$ws.Range("N140").Value = -26258.8

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2505.2727
$ws.Range("I102").Value = 2147.7144
$ws.Range("K102").Value = 2147.7144
$ws.Range("M102").Value = -525.7143999999998
$ws.Range("H122").Value = 3937.2727
$ws.Range("I122").Value = 3237.1428
$ws.Range("K122").Value = 9711.428400000001
$ws.Range("M122").Value = -7261.428400000001
$ws.Range("H126").Value = 3203.2917
$ws.Range("I126").Value = 2391.4443
$ws.Range("J126").Value = 3690.4
$ws.Range("K126").Value = 7174.3329
$ws.Range("L126").Value = 11071.2
$ws.Range("M126").Value = -4704.3329
$ws.Range("N126").Value = -16011.2

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5242.7856
$ws.Range("I7").Value = 5283.25
$ws.Range("K7").Value = 5283.25
$ws.Range("M7").Value = -5171.25
$ws.Range("H22").Value = 6900.1665
$ws.Range("I22").Value = 8200.200000000001
$ws.Range("J22").Value = 400
$ws.Range("K22").Value = 8200.200000000001
$ws.Range("L22").Value = 400
$ws.Range("M22").Value = -7905.200000000001
$ws.Range("N22").Value = -990
$ws.Range("H27").Value = 6900.1665
$ws.Range("I27").Value = 8200.200000000001
$ws.Range("J27").Value = 400
$ws.Range("K27").Value = 8200.200000000001
$ws.Range("L27").Value = 400
$ws.Range("M27").Value = -8093.200000000001
$ws.Range("N27").Value = -614
$ws.Range("H40").Value = 3585.8096
$ws.Range("I40").Value = 3212.0557
$ws.Range("K40").Value = 3212.0557
$ws.Range("M40").Value = -3076.0557
$ws.Range("H46").Value = 2343
$ws.Range("I46").Value = 2233.3333
$ws.Range("K46").Value = 2233.3333
$ws.Range("M46").Value = -2045.3333
$ws.Range("H100").Value = 2518
$ws.Range("J100").Value = 2518
$ws.Range("L100").Value = 2518
$ws.Range("N100").Value = -3600
$ws.Range("H122").Value = 3271637.2
$ws.Range("I122").Value = 3271637.2
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 9814911.600000001
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -9812461.600000001
$ws.Range("N122").ClearContents()
$ws.Range("H126").Value = 5242.7856
$ws.Range("I126").Value = 5283.25
$ws.Range("K126").Value = 15849.75
$ws.Range("M126").Value = -13379.75

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H26").Value = 250005000
$ws.Range("J26").Value = 6666.6665
$ws.Range("L26").Value = 6666.6665
$ws.Range("N26").Value = -7252.6665
$ws.Range("H43").Value = 6200
$ws.Range("J43").Value = 6200
$ws.Range("L43").Value = 6200
$ws.Range("N43").Value = -6498
$ws.Range("H122").Value = 1238.0769
$ws.Range("I122").Value = 963.63635
$ws.Range("K122").Value = 2890.90905
$ws.Range("M122").Value = -440.9090500000002
